# "chinh sua yeu cau thay doi" - insert a new header row above the existing
# "begin header" row with the requested week range, bold/underlined/14pt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at row 3 - everything from the old row 3 down
# (old rows 3..12) shifts down by one (new rows 4..13).
$ws.Rows("3:3").Insert()

# Fill in the new row's A3 cell with the week label and its formatting.
$cell = $ws.Range("A3")
$cell.Value = "Week from 14-03 to 20-03"
$cell.Font.Bold = $true
$cell.Font.Italic = $false
$cell.Font.Underline = $true
$cell.Font.Size = 14
$cell.Font.Name = "Calibri"

# Match the selection left behind by the edit (A4, the cell that used to be
# A3 before the insert).
$ws.Range("A4").Select()
